# johnson_pierce.xlsx regen: replace column G ("K") values, computed from
# the new Strike# (TB) based K-count logic (regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals).
#
# The sheet stores raw literal values (no formulas), so we simply overwrite
# each data row's G cell (K) with its freshly-regenerated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for data rows 2..67 (Excel row number -> new value), in order.
$kValues = @(
    2,0,1,2,2,2,2,2,0,0,
    0,1,0,2,1,2,3,2,0,3,
    3,1,1,3,2,1,1,0,1,0,
    2,1,0,2,1,2,0,1,1,1,
    2,1,1,2,0,1,3,1,3,1,
    2,1,2,2,1,0,2,0,2,2,
    2,0,1,2,1,0
)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
